$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Row 2 updates
$ws.Range("G2").Value = 1.96
$ws.Range("H2").Value = 2.75
$ws.Range("N2").Value = 3.4
$ws.Range("O2").Value = 1.33
$ws.Range("R2").Value = 2.5
$ws.Range("S2").Value = 1.5
$ws.Range("V2").Value = 11
$ws.Range("AD2").Value = 8
$ws.Range("AF2").Value = 17

# Row 3 updates
$ws.Range("K3").Value = 8
$ws.Range("L3").Value = 1.36
$ws.Range("M3").Value = 3

# Row 4 updates
$ws.Range("G4").Value = 2.63
$ws.Range("H4").Value = 2.7
$ws.Range("I4").Value = 3.1
$ws.Range("N4").Value = 2.25
$ws.Range("O4").Value = 1.62
$ws.Range("R4").Value = 1.83
$ws.Range("S4").Value = 1.83
$ws.Range("T4").Value = 8
$ws.Range("U4").Value = 12
$ws.Range("W4").Value = 26
$ws.Range("AA4").Value = 5.5
$ws.Range("AE4").Value = 15
$ws.Range("AF4").Value = 12
$ws.Range("AH4").Value = 26
$ws.Range("AI4").Value = 41
